$wb = $excel.ActiveWorkbook

# --- Update the "Hoja1" sheet: conversion rates text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$nuevoTexto = @"
Conversión del día 💰
✅ Dólar paralelo: 68

Binance
✅ 1000 Bs = 1.79 = 6589.62 pesos
✅ 6589.62 pesos = 1.78 = 940.38 Bs

Promedio competencia
✅ Tasa pesos: 20
✅ Tasa Bs: 20
✅ % Ganancia: 20%
"@
$ws1.Range("A1").Value = $nuevoTexto

# --- Update the "tasas" sheet: updated rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 558
$ws2.Range("O10").Value = 3677.01
$ws2.Range("N12").Value = 3699.94
$ws2.Range("O12").Value = 528.002
